$p = $ppt.ActivePresentation

# 1. Footer "date" placeholder text: 10/16/2014 -> 10/17/2014
#    The placeholder lives on the slide master and is repeated on every
#    slide layout, so walk both collections and update each one.
$m = $p.SlideMaster

for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = "10/17/2014"
    }
}

$layouts = $m.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $lay = $layouts.Item($i)
    $shapes = $lay.Shapes
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "10/17/2014"
        }
    }
}

# 2. Slide 3 heading textbox: "Chassis" -> "Vehicle "
$s3 = $p.Slides.Item(3)
for ($i = 1; $i -le $s3.Shapes.Count; $i++) {
    $sh = $s3.Shapes.Item($i)
    if ($sh.Name -eq "TextBox 123") {
        $sh.TextFrame.TextRange.Text = "Vehicle "
    }
}
